$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.720.88"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "1.847.83"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4304"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07335"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8766"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.820.51"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.530"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06945"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008988"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("D21").Value = "27.746.05"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.980"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E23").Value = "  -2.43%  "

$ws.Range("D24").Value = "2.054.07"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("E25").Value = "  -3.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.265"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "119.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.878"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08901"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7557"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.554"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.969"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.87%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01936"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5094"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1663"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.581"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.387"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06542"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4670"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.632"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.55"
$ws.Range("D51").Style = "Normal"
